$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 108284574
$ws.Range("Q2").Value = 445714.1228795081
$ws.Range("R2").Value = 7066173.645204219

# Row 3
$ws.Range("A3").Value = 108284567
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("Q3").Value = 445434.438710931
$ws.Range("R3").Value = 7065991.394450839
$ws.Range("AC3").Value = "ringhack"

# Row 4
$ws.Range("A4").Value = 108284581
$ws.Range("Q4").Value = 445796.1410013655
$ws.Range("R4").Value = 7066400.520577709

# Row 5
$ws.Range("A5").Value = 108284571
$ws.Range("M5").Value = "gammalt bo"
$ws.Range("N5").Value = ""
$ws.Range("Q5").Value = 445647.4287634333
$ws.Range("R5").Value = 7066063.9690001
$ws.Range("AC5").ClearContents()

# Row 6
$ws.Range("A6").Value = 108284577
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("Q6").Value = 445809.6948888918
$ws.Range("R6").Value = 7066181.755240711
$ws.Range("AC6").Value = "ringhack"

# Row 7
$ws.Range("A7").Value = 108284582
$ws.Range("Q7").Value = 445826.6784053955
$ws.Range("R7").Value = 7066421.272463826

# Row 8
$ws.Range("A8").Value = 108284569
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("Q8").Value = 445629.6519997923
$ws.Range("R8").Value = 7065959.653755097
$ws.Range("AC8").Value = "ringhack gamla"

# Row 9
$ws.Range("A9").Value = 108284575
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("Q9").Value = 445771.9083703306
$ws.Range("R9").Value = 7066178.414627753
$ws.Range("AC9").Value = "ringhack"

# Row 10
$ws.Range("A10").Value = 108284572
$ws.Range("Q10").Value = 445661.8261003256
$ws.Range("R10").Value = 7066100.958404644

# Row 11
$ws.Range("A11").Value = 108284568
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = ""
$ws.Range("Q11").Value = 445617.6319669108
$ws.Range("R11").Value = 7065957.644219733
$ws.Range("AC11").Value = "ringhack gamla"

# Row 12
$ws.Range("A12").Value = 108284573
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("Q12").Value = 445713.419994569
$ws.Range("R12").Value = 7066158.58546807
$ws.Range("AC12").Value = "ringhack"

# Row 13
$ws.Range("A13").Value = 108284570
$ws.Range("Q13").Value = 445645.5316310733
$ws.Range("R13").Value = 7066056.908801682

# Row 14
$ws.Range("A14").Value = 108284576
$ws.Range("Q14").Value = 445808.7769726648
$ws.Range("R14").Value = 7066179.997854604

# Row 15
$ws.Range("A15").Value = 108284580
$ws.Range("Q15").Value = 445814.3812344022
$ws.Range("R15").Value = 7066377.157988134

# Row 16
$ws.Range("A16").Value = 108284579
$ws.Range("Q16").Value = 445845.9283687233
$ws.Range("R16").Value = 7066327.414423619

# Row 17
$ws.Range("A17").Value = 108284651
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = ""
$ws.Range("Q17").Value = 445344.3129855981
$ws.Range("R17").Value = 7066093.141319267
$ws.Range("AC17").Value = "ringhack gamla"
